$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers must be forced to text
# format first, matching the original inlineStr (text) cell type in the sheet.
$ws.Range("D2").Value = "69.021.63"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "3.747.31"
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.17"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.91"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").Value = "3.745.35"
$ws.Range("E7").Value = "  +1.29%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.540"
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("E11").Value = "  +3.02%  "
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.11"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("E14").Value = "  +0.90%  "
$ws.Range("D15").Value = "4.372.99"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").Value = "3.748.65"
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("D17").Value = "69.037.55"
$ws.Range("E17").Value = "  +1.23%  "
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("E19").Value = "  -1.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.21"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.00"
$ws.Range("E21").Value = "  +19.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "492.77"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("E24").Value = "  +6.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.88"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.32"
$ws.Range("E26").Value = "  +0.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.31"
$ws.Range("E27").Value = "  +0.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.14"
$ws.Range("E28").Value = "  +0.85%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.99"
$ws.Range("E30").Value = "  +2.89%  "
$ws.Range("E31").Value = "  +4.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.08"
$ws.Range("E32").Value = "  +2.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.64"
$ws.Range("E33").Value = "  +0.75%  "
$ws.Range("D34").Value = "3.892.28"
$ws.Range("E34").Value = "  +1.26%  "
$ws.Range("E35").Value = "  +0.62%  "
$ws.Range("D36").Value = "3.682.28"
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("E38").Value = "  +1.27%  "
$ws.Range("E39").Value = "  +1.91%  "
$ws.Range("E40").Value = "  +1.77%  "
$ws.Range("E41").Value = "  +1.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.97"
$ws.Range("E42").Value = "  +4.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "433.02"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.59"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.00"
$ws.Range("E45").Value = "  +1.91%  "
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.58"
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.45"
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").Value = "2.793.12"
$ws.Range("E50").Value = "  +1.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0353"
$ws.Range("E51").Value = "  +0.90%  "
